$wb = $excel.ActiveWorkbook

# --- Sheet "Admin": update D2 (Key value) ---
$wsAdmin = $wb.Worksheets.Item("Admin")
$wsAdmin.Range("D2").Value = "7064986A"

# --- Sheet "Jira": update existing rows and append a new row ---
$wsJira = $wb.Worksheets.Item("Jira")
$wsJira.Range("B2").Value = "2098D84319E1DC7FC1B341BAEF3278E0"
$wsJira.Range("A3").Value = "Leave_Accept"
$wsJira.Range("A4").Value = "PersonalDetails"
# Give the new, otherwise-empty B4 cell the same (default) style as B3
# so it is materialized in the sheet just like the existing blank cell.
$wsJira.Range("B4").Style = $wsJira.Range("B3").Style
